# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.406.38'
$ws.Range("E2").Value = '  +3.18%  '
$ws.Range("D3").Value = '1.587.56'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("E4").Value = '  +1.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.12%  '
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  +1.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.42'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.74%  '
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0599'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0887'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.13%  '
$ws.Range("D12").Value = '1.814.89'
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("D13").Value = '1.575.43'
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.74'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("D16").Value = '28.420.38'
$ws.Range("E16").Value = '  +3.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.49'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.81%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0706'
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("E28").Value = '  -0.75%  '
$ws.Range("E29").Value = '  +1.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").Value = '1.398.52'
$ws.Range("E34").Value = '  -3.50%  '
$ws.Range("E35").Value = '  -1.07%  '
$ws.Range("E36").Value = '  -8.80%  '
$ws.Range("E37").Value = '  +1.31%  '
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("E39").Value = '  +8.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.540'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("E44").Value = '  -2.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.980'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").Value = '1.725.99'
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.13'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("E50").Value = '  +13.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0519'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.22%  '
